# 19/12/2025: Update the list
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data for rows 2..11 (row 1 is the header and stays unchanged):
# Job ID | Company | Job Title | Candidate | Status
$data = @(
    @(650, "Dash0",             "Enterprise Account Executive - Nordics / Stockholm", "[deleted]",    "CV Sent"),
    @(721, "Rox",                "ENT AE East",                                         "Graham Hill",  "4th Interview"),
    @(721, "Rox",                "ENT AE East",                                         "Graham Hill",  "4th Interview"),
    @(721, "Rox",                "ENT AE East",                                         "Graham Hill",  "4th Interview"),
    @(721, "Rox",                "ENT AE East",                                         "Graham Hill",  "4th Interview"),
    @(766, "Cogent Security",    "Enterprise Account Executive (US)",                   "Daniel Smith", "2nd Interview"),
    @(796, "Redwood Software",   "Enterprise AE UK (Finance Automation)",               "Alex  Kennedy","1st Interview"),
    @(796, "Redwood Software",   "Enterprise AE UK (Finance Automation)",               "Steven  Jones","CV Sent"),
    @(796, "Redwood Software",   "Enterprise AE UK (Finance Automation)",               "Tom Harries",  "1st Interview"),
    @(813, "Laurel",             "Enterprise Account Executive UK x4",                  "Tom Harries",  "1st Interview")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}
